$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data set (scenario name, probability) for rows 2-26
$data = @(
    @(12345, 0.1314),
    @(2345, 0.0493),
    @(1345, 0.0131),
    @(1245, 0.0493),
    @(1235, 0.0131),
    @(1234, 0.0493),
    @(123, 0.0127),
    @(124, 0.0591),
    @(125, 0.0127),
    @(134, 0.0127),
    @(135, 0.0029),
    @(145, 0.0127),
    @(234, 0.0591),
    @(235, 0.0127),
    @(245, 0.0591),
    @(345, 0.0127),
    @(13, 0.0057),
    @(14, 0.0338),
    @(15, 0.0057),
    @(23, 0.0338),
    @(24, 0.266),
    @(25, 0.0338),
    @(34, 0.0338),
    @(35, 0.0057)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}

# Row 26: scenario "F" with probability
$ws.Cells.Item(26, 1).Value = "F"
$ws.Cells.Item(26, 2).Value = 0.0201

# Row 27: blank row (both cells blank but styled)
$ws.Cells.Item(27, 1).Value = ""
$ws.Cells.Item(27, 2).Value = ""

# Row 28: total check value
$ws.Cells.Item(28, 2).Value = 1.0003

# Apply style index 1 (s="1") to column A rows 2-28, and style index 2 (s="2") to column B rows 2-28
$ws.Range("A2:A28").Style = "Normal"
$ws.Range("B2:B28").Style = "Normal"

# Resize table
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B28"))

# Update sheet view selection
$ws.Range("A30").Select()

# Column widths
$ws.Columns.Item(1).ColumnWidth = 16.3046875
$ws.Columns.Item(2).ColumnWidth = 11.23046875
